# Daily attendance processing - 2025-12-30 08:41:18
# Swap the order of "dnasr281@gmail.com" and "System" in the "Recorded By"
# column (G) wherever the combined string "dnasr281@gmail.com, System"
# appears, turning it into "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
